$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Dinesh Chandimal'
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 'LBW'
$ws.Range("E2").Value = ' Josh Hazlewood'
$ws.Range("J2").Value = 'David Warner'
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 5
$ws.Range("N2").Value = ' Dushmantha Chameera'
$ws.Range("A3").Value = 'Pathum Nissanka'
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 'LBW'
$ws.Range("E3").Value = ' Mitchell Starc'
$ws.Range("J3").Value = 'Aaron Finch(C)'
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = 'Bowled'
$ws.Range("N3").Value = ' Nuwan Pradeep'
$ws.Range("A4").Value = 'Charith Asalanka'
$ws.Range("D4").Value = 'LBW'
$ws.Range("E4").Value = ' Josh Hazlewood'
$ws.Range("J4").Value = 'Mitchell Marsh'
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 7
$ws.Range("N4").Value = ' Dushmantha Chameera'
$ws.Range("A5").Value = 'Dhananjaya de Silva'
$ws.Range("B5").Value = 16
$ws.Range("C5").Value = 8
$ws.Range("E5").Value = ' Pat Cummins'
$ws.Range("J5").Value = 'Steve Smith'
$ws.Range("K5").Value = 8
$ws.Range("M5").Value = 'Caught'
$ws.Range("N5").Value = ' Chamika Karunarathne'
$ws.Range("A6").Value = 'Bhanuka Rajapakse'
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = ' Pat Cummins'
$ws.Range("J6").Value = 'Glenn Maxwell'
$ws.Range("K6").Value = 14
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 'Caught'
$ws.Range("N6").Value = ' Wanindu Hasaranga'
$ws.Range("A7").Value = 'Dasun Shanka(C)'
$ws.Range("B7").Value = 64
$ws.Range("C7").Value = 28
$ws.Range("D7").Value = 'LBW'
$ws.Range("E7").Value = ' Mitchell Starc'
$ws.Range("J7").Value = 'Matthew Wade'
$ws.Range("K7").Value = 35
$ws.Range("L7").Value = 13
$ws.Range("M7").Value = 'Bowled'
$ws.Range("N7").Value = ' Dushmantha Chameera'
$ws.Range("A8").Value = 'Wanindu Hasaranga'
$ws.Range("B8").Value = 31
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = 'Caught'
$ws.Range("E8").Value = ' Adam Zampa'
$ws.Range("J8").Value = 'Marcus Stionis'
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 'LBW'
$ws.Range("N8").Value = ' Wanindu Hasaranga'
$ws.Range("A9").Value = 'Chamika Karunarathne'
$ws.Range("B9").Value = 1
$ws.Range("D9").Value = 'NOT OUT'
$ws.Range("E9").Value = ' '
$ws.Range("J9").Value = 'Pat Cummins'
$ws.Range("K9").Value = 26
$ws.Range("M9").Value = 'LBW'
$ws.Range("N9").Value = ' Maheesh Theekshana'
$ws.Range("A10").Value = 'Dushmantha Chameera'
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 'LBW'
$ws.Range("E10").Value = ' Pat Cummins'
$ws.Range("J10").Value = 'Mitchell Starc'
$ws.Range("K10").Value = 32
$ws.Range("L10").Value = 9
$ws.Range("N10").Value = ' Nuwan Pradeep'
$ws.Range("A11").Value = 'Maheesh Theekshana'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 1
$ws.Range("E11").Value = ' Pat Cummins'
$ws.Range("J11").Value = 'Adam Zampa'
$ws.Range("K11").Value = 32
$ws.Range("L11").Value = 9
$ws.Range("M11").Value = '* NOT OUT'
$ws.Range("A12").Value = 'Nuwan Pradeep'
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = 2
$ws.Range("E12").Value = ' Pat Cummins'
$ws.Range("J12").Value = 'Josh Hazlewood'
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 'NOT OUT'
$ws.Range("N12").Value = ' '
$ws.Range("A16").Value = 166
$ws.Range("C16").Value = '''13.4'
$ws.Range("D16").Value = 82
$ws.Range("J16").Value = 170
$ws.Range("K16").Value = 9
$ws.Range("L16").Value = '''10.4'
$ws.Range("M16").Value = 64
$ws.Range("A21").Value = 'Marcus Stionis'
$ws.Range("J21").Value = 'Maheesh Theekshana'
$ws.Range("L21").Value = 34
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 17
$ws.Range("A22").Value = 'Josh Hazlewood'
$ws.Range("J22").Value = 'Dushmantha Chameera'
$ws.Range("L22").Value = 30
$ws.Range("M22").Value = 3
$ws.Range("N22").Value = 15
$ws.Range("A23").Value = 'Adam Zampa'
$ws.Range("C23").Value = 35
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 11.67
$ws.Range("J23").Value = 'Chamika Karunarathne'
$ws.Range("L23").Value = 31
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 15.5
$ws.Range("A24").Value = 'Mitchell Starc'
$ws.Range("C24").Value = 31
$ws.Range("E24").Value = 10.33
$ws.Range("J24").Value = 'Wanindu Hasaranga'
$ws.Range("L24").Value = 34
$ws.Range("N24").Value = 17
$ws.Range("A25").Value = 'Pat Cummins'
$ws.Range("B25").Value = '''2.4'
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 10.42
$ws.Range("J25").Value = 'Nuwan Pradeep'
$ws.Range("K25").Value = '''2.4'
$ws.Range("L25").Value = 41
$ws.Range("N25").Value = 17.08

Write-Host "Applied all updates"
